$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the applications table (A1:I11) by the "id" column ascending, keeping row 1 as header.
$sortRange = $ws.Range("A1:I11")
$sortKey = $ws.Range("A2:A11")
$sortRange.Sort($sortKey, 1, $null, $null, 1, 0, 1, 1)

# After sorting, row 4 now holds id=2 ("Mommy" test record) and rows 7-11 hold ids 5-9
# (leftover test/demo rows). Remove them so only the header plus ids 0, 1, 3, 4 remain,
# matching the trimmed-down "app form" table.
$ws.Range("A7:A11").EntireRow.Delete()
$ws.Range("A4").EntireRow.Delete()
